$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2431.5862
$ws.Range("I15").Value = 2431.5862
$ws.Range("K15").Value = 7294.758600000001
$ws.Range("M15").Value = -7125.758600000001
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("H32").Value = 655.1667
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 673.55554
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 673.55554
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1325.55554
$ws.Range("H33").Value = 543.4815
$ws.Range("I33").Value = 415.78946
$ws.Range("J33").Value = 846.75
$ws.Range("K33").Value = 415.78946
$ws.Range("L33").Value = 846.75
$ws.Range("M33").Value = -186.78946
$ws.Range("N33").Value = -1304.75
$ws.Range("H69").Value = 4913
$ws.Range("I69").Value = 4913
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 14739
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -13865
$ws.Range("H72").Value = 4913
$ws.Range("I72").Value = 4913
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 44217
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -39849
$ws.Range("H132").Value = 3706248.2
$ws.Range("I132").Value = 4257611
$ws.Range("J132").Value = 4243.2856
$ws.Range("K132").Value = 12772833
$ws.Range("L132").Value = 12729.8568
$ws.Range("M132").Value = -12770303
$ws.Range("N132").Value = -17789.8568
$ws.Range("H135").Value = 1236.7142
$ws.Range("I135").Value = 942.8333
$ws.Range("K135").Value = 8485.4997
$ws.Range("M135").Value = -5950.4997
$ws.Range("H137").Value = 2502968
$ws.Range("I137").Value = 4351483.5
$ws.Range("K137").Value = 13054450.5
$ws.Range("M137").Value = -13051900.5
$ws.Range("H138").Value = 2472.9673
$ws.Range("J138").Value = 4313.9766
$ws.Range("L138").Value = 12941.9298
$ws.Range("N138").Value = -23221.9298
$ws.Range("H141").Value = 174640.98
$ws.Range("I141").Value = 923.4528
$ws.Range("J141").Value = 1197644.2
$ws.Range("K141").Value = 2770.3584
$ws.Range("L141").Value = 3592932.6
$ws.Range("M141").Value = 2409.6416
$ws.Range("N141").Value = -3603292.6
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3598.9
$ws.Range("I32").Value = 2963.9302
$ws.Range("K32").Value = 2963.9302
$ws.Range("M32").Value = -2676.9302
$ws.Range("H74").Value = 900.8
$ws.Range("I74").Value = 778.7222
$ws.Range("J74").Value = 1999.5
$ws.Range("K74").Value = 778.7222
$ws.Range("L74").Value = 1999.5
$ws.Range("M74").Value = 95.27779999999996
$ws.Range("N74").Value = -3747.5
$ws.Range("H77").Value = 900.8
$ws.Range("I77").Value = 778.7222
$ws.Range("J77").Value = 1999.5
$ws.Range("K77").Value = 3893.611
$ws.Range("L77").Value = 9997.5
$ws.Range("M77").Value = 474.3889999999997
$ws.Range("N77").Value = -18733.5
$ws.Range("H141").Value = 29666.666
$ws.Range("J141").Value = 29666.666
$ws.Range("L141").Value = 29666.666
$ws.Range("N141").Value = -40026.666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5595
$ws.Range("I134").Value = 2800
$ws.Range("K134").Value = 8400
$ws.Range("M134").Value = -5865

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1788356.5
$ws.Range("I31").Value = 2633223.8
$ws.Range("J31").Value = 4747.8335
$ws.Range("K31").Value = 2633223.8
$ws.Range("L31").Value = 4747.8335
$ws.Range("M31").Value = -2632928.8
$ws.Range("N31").Value = -5337.8335
$ws.Range("H34").Value = 1788356.5
$ws.Range("I34").Value = 2633223.8
$ws.Range("J34").Value = 4747.8335
$ws.Range("K34").Value = 2633223.8
$ws.Range("L34").Value = 4747.8335
$ws.Range("M34").Value = -2633021.8
$ws.Range("N34").Value = -5151.8335
$ws.Range("H134").Value = 2492.9565
$ws.Range("I134").Value = 1326.5625
$ws.Range("K134").Value = 3979.6875
$ws.Range("M134").Value = -1444.6875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 8232.666999999999
$ws.Range("I34").Value = 215
$ws.Range("J34").Value = 13577.777
$ws.Range("K34").Value = 645
$ws.Range("L34").Value = 40733.331
$ws.Range("M34").Value = -561
$ws.Range("N34").Value = -40901.331
$ws.Range("H40").Value = 192.85715
$ws.Range("I40").Value = 100
$ws.Range("J40").Value = 230
$ws.Range("K40").Value = 400
$ws.Range("L40").Value = 920
$ws.Range("M40").Value = -331
$ws.Range("N40").Value = -1058
$ws.Range("H60").Value = 14671.429
$ws.Range("I60").Value = 450
$ws.Range("J60").Value = 100000
$ws.Range("K60").Value = 1350
$ws.Range("L60").Value = 300000
$ws.Range("M60").Value = -1099
$ws.Range("N60").Value = -300502
$ws.Range("H132").Value = 2701
$ws.Range("I132").Value = 2034.6666
$ws.Range("J132").Value = 3034.1667
$ws.Range("K132").Value = 18311.9994
$ws.Range("L132").Value = 27307.5003
$ws.Range("M132").Value = -15781.9994
$ws.Range("N132").Value = -32367.5003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 32500
$ws.Range("J118").Value = 32500
$ws.Range("L118").Value = 32500
$ws.Range("N118").Value = -35814
$ws.Range("H132").Value = 2714.1875
$ws.Range("I132").Value = 2402.25
$ws.Range("J132").Value = 3650
$ws.Range("K132").Value = 7206.75
$ws.Range("L132").Value = 10950
$ws.Range("M132").Value = -4676.75
$ws.Range("N132").Value = -16010
$ws.Range("H135").Value = 29535
$ws.Range("J135").Value = 29535
$ws.Range("L135").Value = 29535
$ws.Range("N135").Value = -39675

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1931.1666
$ws.Range("I7").Value = 1130.4445
$ws.Range("J7").Value = 4333.3335
$ws.Range("K7").Value = 1130.4445
$ws.Range("L7").Value = 4333.3335
$ws.Range("M7").Value = -1018.4445
$ws.Range("N7").Value = -4557.3335
$ws.Range("H126").Value = 1931.1666
$ws.Range("I126").Value = 1130.4445
$ws.Range("J126").Value = 4333.3335
$ws.Range("K126").Value = 3391.3335
$ws.Range("L126").Value = 13000.0005
$ws.Range("M126").Value = -921.3335000000002
$ws.Range("N126").Value = -17940.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 29000
$ws.Range("J76").Value = 29000
$ws.Range("L76").Value = 29000
$ws.Range("N76").Value = -29630
$ws.Range("H79").Value = 29000
$ws.Range("J79").Value = 29000
$ws.Range("L79").Value = 29000
$ws.Range("N79").Value = -31184
$ws.Range("H126").Value = 3334670.8
$ws.Range("I126").Value = 744.96295
$ws.Range("K126").Value = 2234.88885
$ws.Range("M126").Value = 235.1111500000002

